# Applies the "Sync attendance_reports..." edit to the session analysis sheet:
#  1. Widen the Subject column (C).
#  2. Rename the subject "GENERAL SURGERY" -> "SURGERY SEMINAR/SLIDE" for every
#     session row of every B1 group.
#  3. Reorder the "Recorded By" text for the first three sessions of each group.
#  4. Flip session #16 of every group from Recorded -> Not Recorded (clear the
#     recorder, zero the attendance count, change status + row colour).
#  5. Refresh the global summary box (Recorded/Missing/Coverage/Average).
#  6. Refresh the per-group statistics table (Recorded/Missing/Coverage/Average).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Column C width (target stored width = 23) ---------------------------
$ws.Columns("C").ColumnWidth = 22.166666666666668

# --- group definitions --------------------------------------------------
# Start = first session row (session 1) for the group in the A:I block
# Total = class size (denominator of the H "x/y" attendance fractions)
# StatsRow = the row of this group inside the K:S "Group Statistics" table
# O/P/R/S = the new Recorded / Missing / Coverage% / Avg Attendance% values
$groups = @(
    @{ Start=2;   Total=27; StatsRow=15; O=4; P=2; R="18.2%"; S="82.4%" },
    @{ Start=24;  Total=31; StatsRow=16; O=3; P=2; R="13.6%"; S="72.0%" },
    @{ Start=46;  Total=19; StatsRow=17; O=3; P=2; R="13.6%"; S="57.9%" },
    @{ Start=68;  Total=21; StatsRow=18; O=3; P=2; R="13.6%"; S="90.5%" },
    @{ Start=90;  Total=31; StatsRow=19; O=4; P=2; R="18.2%"; S="96.0%" },
    @{ Start=112; Total=28; StatsRow=20; O=4; P=2; R="18.2%"; S="92.9%" },
    @{ Start=134; Total=29; StatsRow=21; O=4; P=2; R="18.2%"; S="92.2%" },
    @{ Start=156; Total=33; StatsRow=22; O=4; P=2; R="18.2%"; S="89.4%" },
    @{ Start=178; Total=30; StatsRow=23; O=4; P=2; R="18.2%"; S="84.2%" },
    @{ Start=200; Total=27; StatsRow=24; O=3; P=2; R="13.6%"; S="65.4%" },
    @{ Start=222; Total=29; StatsRow=25; O=3; P=2; R="13.6%"; S="73.6%" },
    @{ Start=244; Total=29; StatsRow=26; O=3; P=2; R="13.6%"; S="63.2%" }
)

$pink = 12695295   # RGB(255,182,193) - "Not Recorded" row colour

foreach ($g in $groups) {
    $start = $g.Start

    # --- 2. Rename the subject for all 22 sessions of this group -------
    for ($i = 0; $i -lt 22; $i++) {
        $row = $start + $i
        if ($ws.Cells.Item($row, 3).Value -eq "GENERAL SURGERY") {
            $ws.Cells.Item($row, 3).Value = "SURGERY SEMINAR/SLIDE"
        }
    }

    # --- 3. Reorder "Recorded By" for sessions 1-3 ----------------------
    for ($i = 0; $i -lt 3; $i++) {
        $row = $start + $i
        if ($ws.Cells.Item($row, 7).Value -eq "dnasr281@gmail.com, System") {
            $ws.Cells.Item($row, 7).Value = "System, dnasr281@gmail.com"
        }
    }

    # --- 4. Session #16 becomes "Not Recorded" --------------------------
    $row16 = $start + 15
    $ws.Range($ws.Cells.Item($row16, 1), $ws.Cells.Item($row16, 9)).Interior.Color = $pink
    $ws.Cells.Item($row16, 7).Value = ""
    $ws.Cells.Item($row16, 8).Value = "0/" + $g.Total
    $ws.Cells.Item($row16, 9).Value = "Not Recorded"

    # --- 6. Per-group statistics table -----------------------------------
    $sr = $g.StatsRow
    $ws.Cells.Item($sr, 15).Value = $g.O   # Column O - Recorded
    $ws.Cells.Item($sr, 16).Value = $g.P   # Column P - Missing
    $ws.Cells.Item($sr, 18).Value = $g.R   # Column R - Coverage %
    $ws.Cells.Item($sr, 19).Value = $g.S   # Column S - Avg Attendance %
}

# --- 5. Global summary box -------------------------------------------------
$ws.Cells.Item(6, 12).Value = 42       # Recorded Sessions
$ws.Cells.Item(7, 12).Value = 24       # Missing Sessions
$ws.Cells.Item(9, 12).Value = "15.9%"  # Coverage %
$ws.Cells.Item(10, 12).Value = "81.3%" # Average Attendance %
